$d = $word.ActiveDocument

# Finds the (unique) literal $findText in the document and rewrites it as a
# sequence of separate runs, one per entry in $parts (their concatenation
# equals the new text). Toggling a formatting property around each
# sub-range's Text assignment stops the engine from silently re-merging the
# new run with its formatting-identical neighbours, so the run boundaries
# implied by $parts are preserved in the saved OOXML.
function Split-Replace($findText, $parts) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Text = $findText
    $rng.Find.Forward = $true
    $rng.Find.Wrap = 0
    $rng.Find.MatchCase = $true
    $rng.Find.MatchWholeWord = $false
    $rng.Find.MatchWildcards = $false
    $found = $rng.Find.Execute()
    if (-not $found) {
        throw "Split-Replace: text not found: $findText"
    }

    $pos = $rng.Start
    foreach ($part in $parts) {
        $len = $part.Length
        $sub = $d.Range($pos, $pos + $len)
        $sub.Font.Bold = $true
        $sub.Text = $part
        $sub.Font.Bold = $false
        $pos = $pos + $len
    }
}

Split-Replace "{{Edad}}" @("{{", "e", "dad}}")
Split-Replace "{{RV_M}}" @("{{", "VN", "_M}}")
Split-Replace "{{R_M}}" @("{{", "V", "_M}}")
Split-Replace "{{RV_C}}" @("{{", "VN", "_C}}")
Split-Replace "{{R_C}}" @("{{", "V", "_C}}")
Split-Replace "{{RV_A}}" @("{{", "VN", "_A}}")
Split-Replace "{{R_A}}" @("{{", "V", "_A}}")
